# Remove the prod_e / prod_l prediction columns (H:I) that were added to
# the original prediction_data sheet, restoring it to the "data only"
# layout used for prediction (no prod_e / prod_l labels or values).
#
# H1/I1 keep their header style but lose their shared-string text.
# H2:I92 (the data rows) are cleared entirely.
# Selection moves to reflect the now-empty H1:I92 block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the header labels in H1:I1 (style is preserved, text removed).
$ws.Range("H1:I1").ClearContents()

# Clear out the prod_e / prod_l data values for all data rows.
$ws.Range("H2:I92").ClearContents()

# Reflect the new selection/viewport over the cleared block.
$ws.Range("H1:I92").Select()
